$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet ALC - row 138 : refreshed market data
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4999.75
$ws.Range("J138").Value = 4999.75
$ws.Range("L138").Value = 14999.25
$ws.Range("N138").Value = -25279.25

# ---------------------------------------------------------------
# Sheet CRP - refreshed / cleared market data
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 58: data now all zero, profit columns no longer populated
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58:N58").ClearContents()

# Rows with no market data fetched this run: clear H:N entirely
foreach ($r in 129,130,131,132,133,134,135,137,138,139,140,141) {
    $ws.Range("H" + $r + ":N" + $r).ClearContents()
}

# Row 136: data now all zero, profit columns no longer populated
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136:N136").ClearContents()

# ---------------------------------------------------------------
# Sheet CUL - refreshed / cleared market data
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H11").Value = 190
$ws.Range("I11").Value = 190
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 570
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -430
$ws.Range("N11").ClearContents()

$ws.Range("H68").Value = 4999.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4999.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 14998.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -16620.5

$ws.Range("H71").Value = 4999.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4999.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 44995.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -53107.5

# Rows with no market data fetched this run: clear H:N entirely
# (row 135 is left untouched)
foreach ($r in 120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,136,137,138,139,140,141) {
    $ws.Range("H" + $r + ":N" + $r).ClearContents()
}

# ---------------------------------------------------------------
# Sheet GSM - row 27 : refreshed market data
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
